$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force the value to be stored as text even when it looks like a number
    # (e.g. "212.23"), then reset the style so no extra "Text" number format
    # is left behind on the cell.
    $ws.Range($range).Value = "'" + $value
    $ws.Range($range).Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.253.69"
$ws.Range("E2").Value = "  +0.40%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.588.88"
$ws.Range("E3").Value = "  +0.73%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.13%  "

# Row 5 - BNB
Set-TextValue "D5" "212.23"
$ws.Range("E5").Value = "  +1.44%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.74%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.12%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +0.14%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.0607"
$ws.Range("E9").Value = "  -0.28%  "

# Row 10 - Solana
Set-TextValue "D10" "19.33"
$ws.Range("E10").Value = "  -1.06%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0850"
$ws.Range("E11").Value = "  +0.74%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.812.49"

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.581.85"
$ws.Range("E13").Value = "  +0.12%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.41%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +1.17%  "

# Row 16 - Litecoin
Set-TextValue "D16" "64.37"
$ws.Range("E16").Value = "  -0.11%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.241.20"

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  -0.35%  "

# Row 19 - Chainlink
Set-TextValue "D19" "7.41"
$ws.Range("E19").Value = "  +2.16%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "213.37"
$ws.Range("E20").Value = "  +2.87%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.13%  "

# Row 22 - Uniswap
Set-TextValue "D22" "4.28"
$ws.Range("E22").Value = "  +0.97%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  +1.31%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -2.53%  "

# Row 25 - Monero
Set-TextValue "D25" "143.93"
$ws.Range("E25").Value = "  +0.02%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  +1.12%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -0.35%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -1.68%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +1.14%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -0.24%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -1.44%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.332.44"
$ws.Range("E34").Value = "  +4.21%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  -0.73%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  -0.63%  "

# Row 37 - ImmutableX
Set-TextValue "D37" "0.587"
$ws.Range("E37").Value = "  -3.32%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +0.59%  "

# Row 39 - ARBITRUM
Set-TextValue "D39" "0.817"
$ws.Range("E39").Value = "  +0.01%  "

# Row 40 - WEMIXToken
$ws.Range("E40").Value = "  -8.11%  "

# Row 41 - FraxShare
Set-TextValue "D41" "5.72"

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  -0.12%  "

# Row 43 - MXToken
$ws.Range("E43").Value = "  +0.09%  "

# Row 44 - TrustWalletToken
$ws.Range("E44").Value = "  +0.31%  "

# Row 45 - Aave
Set-TextValue "D45" "61.80"
$ws.Range("E45").Value = "  -1.03%  "

# Row 46 - RocketPoolETH
$ws.Range("D46").Value = "1.724.14"
$ws.Range("E46").Value = "  +0.70%  "

# Row 47 - Quant
Set-TextValue "D47" "85.87"
$ws.Range("E47").Value = "  -3.48%  "

# Row 48 - RenderToken
Set-TextValue "D48" "1.48"
$ws.Range("E48").Value = "  -3.41%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  -0.72%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  -3.35%  "

# Row 51 - USDD
$ws.Range("E51").Value = "  -0.22%  "
